# Auto-generated edit script: fixes 22 rows in "Stock Report" sheet where the
# damage-detail columns (W:AB) had been collapsed into a single garbled
# Ruby-style array literal stored in column B. This spreads that data back
# out across columns W (Damage Area), X (Damage Component), Y (Remark),
# and Z/AA/AB (blank trailing fields), and blanks out column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stock Report")

# Row 9
$ws.Cells.Item(9, 2).Value = ""
$ws.Cells.Item(9, 2).Style = "Normal"
$ws.Cells.Item(9, 3).Value = ""
$ws.Cells.Item(9, 3).Style = "Normal"
$ws.Cells.Item(9, 4).Value = ""
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = ""
$ws.Cells.Item(9, 5).Style = "Normal"
$ws.Cells.Item(9, 6).Value = ""
$ws.Cells.Item(9, 6).Style = "Normal"
$ws.Cells.Item(9, 7).Value = ""
$ws.Cells.Item(9, 7).Style = "Normal"
$ws.Cells.Item(9, 8).Value = ""
$ws.Cells.Item(9, 8).Style = "Normal"
$ws.Cells.Item(9, 9).Value = ""
$ws.Cells.Item(9, 9).Style = "Normal"
$ws.Cells.Item(9, 10).Value = ""
$ws.Cells.Item(9, 10).Style = "Normal"
$ws.Cells.Item(9, 11).Value = ""
$ws.Cells.Item(9, 11).Style = "Normal"
$ws.Cells.Item(9, 12).Value = ""
$ws.Cells.Item(9, 12).Style = "Normal"
$ws.Cells.Item(9, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(9, 14).Value = ""
$ws.Cells.Item(9, 14).Style = "Normal"
$ws.Cells.Item(9, 15).Value = ""
$ws.Cells.Item(9, 15).Style = "Normal"
$ws.Cells.Item(9, 16).Value = ""
$ws.Cells.Item(9, 16).Style = "Normal"
$ws.Cells.Item(9, 17).Value = ""
$ws.Cells.Item(9, 17).Style = "Normal"
$ws.Cells.Item(9, 18).Value = ""
$ws.Cells.Item(9, 18).Style = "Normal"
$ws.Cells.Item(9, 19).Value = ""
$ws.Cells.Item(9, 19).Style = "Normal"
$ws.Cells.Item(9, 20).Value = ""
$ws.Cells.Item(9, 20).Style = "Normal"
$ws.Cells.Item(9, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(9, 22).Value = ""
$ws.Cells.Item(9, 22).Style = "Normal"
$ws.Cells.Item(9, 23).Value = "DOORS-(D)"
$ws.Cells.Item(9, 24).Value = "Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)"
$ws.Cells.Item(9, 25).Value = "RIGHT DOOR BOTTOM GASKET LOOSE 12'"
$ws.Cells.Item(9, 26).Value = ""
$ws.Cells.Item(9, 26).Style = "Normal"
$ws.Cells.Item(9, 27).Value = ""
$ws.Cells.Item(9, 27).Style = "Normal"
$ws.Cells.Item(9, 28).Value = ""
$ws.Cells.Item(9, 28).Style = "Normal"

# Row 14
$ws.Cells.Item(14, 2).Value = ""
$ws.Cells.Item(14, 2).Style = "Normal"
$ws.Cells.Item(14, 3).Value = ""
$ws.Cells.Item(14, 3).Style = "Normal"
$ws.Cells.Item(14, 4).Value = ""
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = ""
$ws.Cells.Item(14, 5).Style = "Normal"
$ws.Cells.Item(14, 6).Value = ""
$ws.Cells.Item(14, 6).Style = "Normal"
$ws.Cells.Item(14, 7).Value = ""
$ws.Cells.Item(14, 7).Style = "Normal"
$ws.Cells.Item(14, 8).Value = ""
$ws.Cells.Item(14, 8).Style = "Normal"
$ws.Cells.Item(14, 9).Value = ""
$ws.Cells.Item(14, 9).Style = "Normal"
$ws.Cells.Item(14, 10).Value = ""
$ws.Cells.Item(14, 10).Style = "Normal"
$ws.Cells.Item(14, 11).Value = ""
$ws.Cells.Item(14, 11).Style = "Normal"
$ws.Cells.Item(14, 12).Value = ""
$ws.Cells.Item(14, 12).Style = "Normal"
$ws.Cells.Item(14, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(14, 14).Value = ""
$ws.Cells.Item(14, 14).Style = "Normal"
$ws.Cells.Item(14, 15).Value = ""
$ws.Cells.Item(14, 15).Style = "Normal"
$ws.Cells.Item(14, 16).Value = ""
$ws.Cells.Item(14, 16).Style = "Normal"
$ws.Cells.Item(14, 17).Value = ""
$ws.Cells.Item(14, 17).Style = "Normal"
$ws.Cells.Item(14, 18).Value = ""
$ws.Cells.Item(14, 18).Style = "Normal"
$ws.Cells.Item(14, 19).Value = ""
$ws.Cells.Item(14, 19).Style = "Normal"
$ws.Cells.Item(14, 20).Value = ""
$ws.Cells.Item(14, 20).Style = "Normal"
$ws.Cells.Item(14, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(14, 22).Value = ""
$ws.Cells.Item(14, 22).Style = "Normal"
$ws.Cells.Item(14, 23).Value = "PANELS-(PANELS)"
$ws.Cells.Item(14, 24).Value = "Right side panel -(Right side panel )"
$ws.Cells.Item(14, 25).Value = "INTERIOR PANEL DIRTY"
$ws.Cells.Item(14, 26).Value = ""
$ws.Cells.Item(14, 26).Style = "Normal"
$ws.Cells.Item(14, 27).Value = ""
$ws.Cells.Item(14, 27).Style = "Normal"
$ws.Cells.Item(14, 28).Value = ""
$ws.Cells.Item(14, 28).Style = "Normal"

# Row 16
$ws.Cells.Item(16, 2).Value = ""
$ws.Cells.Item(16, 2).Style = "Normal"
$ws.Cells.Item(16, 3).Value = ""
$ws.Cells.Item(16, 3).Style = "Normal"
$ws.Cells.Item(16, 4).Value = ""
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = ""
$ws.Cells.Item(16, 5).Style = "Normal"
$ws.Cells.Item(16, 6).Value = ""
$ws.Cells.Item(16, 6).Style = "Normal"
$ws.Cells.Item(16, 7).Value = ""
$ws.Cells.Item(16, 7).Style = "Normal"
$ws.Cells.Item(16, 8).Value = ""
$ws.Cells.Item(16, 8).Style = "Normal"
$ws.Cells.Item(16, 9).Value = ""
$ws.Cells.Item(16, 9).Style = "Normal"
$ws.Cells.Item(16, 10).Value = ""
$ws.Cells.Item(16, 10).Style = "Normal"
$ws.Cells.Item(16, 11).Value = ""
$ws.Cells.Item(16, 11).Style = "Normal"
$ws.Cells.Item(16, 12).Value = ""
$ws.Cells.Item(16, 12).Style = "Normal"
$ws.Cells.Item(16, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(16, 14).Value = ""
$ws.Cells.Item(16, 14).Style = "Normal"
$ws.Cells.Item(16, 15).Value = ""
$ws.Cells.Item(16, 15).Style = "Normal"
$ws.Cells.Item(16, 16).Value = ""
$ws.Cells.Item(16, 16).Style = "Normal"
$ws.Cells.Item(16, 17).Value = ""
$ws.Cells.Item(16, 17).Style = "Normal"
$ws.Cells.Item(16, 18).Value = ""
$ws.Cells.Item(16, 18).Style = "Normal"
$ws.Cells.Item(16, 19).Value = ""
$ws.Cells.Item(16, 19).Style = "Normal"
$ws.Cells.Item(16, 20).Value = ""
$ws.Cells.Item(16, 20).Style = "Normal"
$ws.Cells.Item(16, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(16, 22).Value = ""
$ws.Cells.Item(16, 22).Style = "Normal"
$ws.Cells.Item(16, 23).Value = "FLOORS-(F)"
$ws.Cells.Item(16, 24).Value = "Threshold plate-(Threshold plate)"
$ws.Cells.Item(16, 25).Value = "FLOOR BOARD NAILS FITTING 8 PES"
$ws.Cells.Item(16, 26).Value = ""
$ws.Cells.Item(16, 26).Style = "Normal"
$ws.Cells.Item(16, 27).Value = ""
$ws.Cells.Item(16, 27).Style = "Normal"
$ws.Cells.Item(16, 28).Value = ""
$ws.Cells.Item(16, 28).Style = "Normal"

# Row 18
$ws.Cells.Item(18, 2).Value = ""
$ws.Cells.Item(18, 2).Style = "Normal"
$ws.Cells.Item(18, 3).Value = ""
$ws.Cells.Item(18, 3).Style = "Normal"
$ws.Cells.Item(18, 4).Value = ""
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = ""
$ws.Cells.Item(18, 5).Style = "Normal"
$ws.Cells.Item(18, 6).Value = ""
$ws.Cells.Item(18, 6).Style = "Normal"
$ws.Cells.Item(18, 7).Value = ""
$ws.Cells.Item(18, 7).Style = "Normal"
$ws.Cells.Item(18, 8).Value = ""
$ws.Cells.Item(18, 8).Style = "Normal"
$ws.Cells.Item(18, 9).Value = ""
$ws.Cells.Item(18, 9).Style = "Normal"
$ws.Cells.Item(18, 10).Value = ""
$ws.Cells.Item(18, 10).Style = "Normal"
$ws.Cells.Item(18, 11).Value = ""
$ws.Cells.Item(18, 11).Style = "Normal"
$ws.Cells.Item(18, 12).Value = ""
$ws.Cells.Item(18, 12).Style = "Normal"
$ws.Cells.Item(18, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(18, 14).Value = ""
$ws.Cells.Item(18, 14).Style = "Normal"
$ws.Cells.Item(18, 15).Value = ""
$ws.Cells.Item(18, 15).Style = "Normal"
$ws.Cells.Item(18, 16).Value = ""
$ws.Cells.Item(18, 16).Style = "Normal"
$ws.Cells.Item(18, 17).Value = ""
$ws.Cells.Item(18, 17).Style = "Normal"
$ws.Cells.Item(18, 18).Value = ""
$ws.Cells.Item(18, 18).Style = "Normal"
$ws.Cells.Item(18, 19).Value = ""
$ws.Cells.Item(18, 19).Style = "Normal"
$ws.Cells.Item(18, 20).Value = ""
$ws.Cells.Item(18, 20).Style = "Normal"
$ws.Cells.Item(18, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(18, 22).Value = ""
$ws.Cells.Item(18, 22).Style = "Normal"
$ws.Cells.Item(18, 23).Value = "PANELS-(PANELS)"
$ws.Cells.Item(18, 24).Value = "Right side panel -(Right side panel )"
$ws.Cells.Item(18, 25).Value = "RIGHT SIDE PANEL PUSH IN 20"" X 40"
$ws.Cells.Item(18, 26).Value = ""
$ws.Cells.Item(18, 26).Style = "Normal"
$ws.Cells.Item(18, 27).Value = ""
$ws.Cells.Item(18, 27).Style = "Normal"
$ws.Cells.Item(18, 28).Value = ""
$ws.Cells.Item(18, 28).Style = "Normal"

# Row 19
$ws.Cells.Item(19, 2).Value = ""
$ws.Cells.Item(19, 2).Style = "Normal"
$ws.Cells.Item(19, 3).Value = ""
$ws.Cells.Item(19, 3).Style = "Normal"
$ws.Cells.Item(19, 4).Value = ""
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = ""
$ws.Cells.Item(19, 5).Style = "Normal"
$ws.Cells.Item(19, 6).Value = ""
$ws.Cells.Item(19, 6).Style = "Normal"
$ws.Cells.Item(19, 7).Value = ""
$ws.Cells.Item(19, 7).Style = "Normal"
$ws.Cells.Item(19, 8).Value = ""
$ws.Cells.Item(19, 8).Style = "Normal"
$ws.Cells.Item(19, 9).Value = ""
$ws.Cells.Item(19, 9).Style = "Normal"
$ws.Cells.Item(19, 10).Value = ""
$ws.Cells.Item(19, 10).Style = "Normal"
$ws.Cells.Item(19, 11).Value = ""
$ws.Cells.Item(19, 11).Style = "Normal"
$ws.Cells.Item(19, 12).Value = ""
$ws.Cells.Item(19, 12).Style = "Normal"
$ws.Cells.Item(19, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(19, 14).Value = ""
$ws.Cells.Item(19, 14).Style = "Normal"
$ws.Cells.Item(19, 15).Value = ""
$ws.Cells.Item(19, 15).Style = "Normal"
$ws.Cells.Item(19, 16).Value = ""
$ws.Cells.Item(19, 16).Style = "Normal"
$ws.Cells.Item(19, 17).Value = ""
$ws.Cells.Item(19, 17).Style = "Normal"
$ws.Cells.Item(19, 18).Value = ""
$ws.Cells.Item(19, 18).Style = "Normal"
$ws.Cells.Item(19, 19).Value = ""
$ws.Cells.Item(19, 19).Style = "Normal"
$ws.Cells.Item(19, 20).Value = ""
$ws.Cells.Item(19, 20).Style = "Normal"
$ws.Cells.Item(19, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(19, 22).Value = ""
$ws.Cells.Item(19, 22).Style = "Normal"
$ws.Cells.Item(19, 23).Value = "PANELS-(PANELS)"
$ws.Cells.Item(19, 24).Value = "Right side panel -(Right side panel )"
$ws.Cells.Item(19, 25).Value = "LEFT SIDE PANEL PUSH IN 20"" 40"""
$ws.Cells.Item(19, 26).Value = ""
$ws.Cells.Item(19, 26).Style = "Normal"
$ws.Cells.Item(19, 27).Value = ""
$ws.Cells.Item(19, 27).Style = "Normal"
$ws.Cells.Item(19, 28).Value = ""
$ws.Cells.Item(19, 28).Style = "Normal"

# Row 22
$ws.Cells.Item(22, 2).Value = ""
$ws.Cells.Item(22, 2).Style = "Normal"
$ws.Cells.Item(22, 3).Value = ""
$ws.Cells.Item(22, 3).Style = "Normal"
$ws.Cells.Item(22, 4).Value = ""
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = ""
$ws.Cells.Item(22, 5).Style = "Normal"
$ws.Cells.Item(22, 6).Value = ""
$ws.Cells.Item(22, 6).Style = "Normal"
$ws.Cells.Item(22, 7).Value = ""
$ws.Cells.Item(22, 7).Style = "Normal"
$ws.Cells.Item(22, 8).Value = ""
$ws.Cells.Item(22, 8).Style = "Normal"
$ws.Cells.Item(22, 9).Value = ""
$ws.Cells.Item(22, 9).Style = "Normal"
$ws.Cells.Item(22, 10).Value = ""
$ws.Cells.Item(22, 10).Style = "Normal"
$ws.Cells.Item(22, 11).Value = ""
$ws.Cells.Item(22, 11).Style = "Normal"
$ws.Cells.Item(22, 12).Value = ""
$ws.Cells.Item(22, 12).Style = "Normal"
$ws.Cells.Item(22, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(22, 14).Value = ""
$ws.Cells.Item(22, 14).Style = "Normal"
$ws.Cells.Item(22, 15).Value = ""
$ws.Cells.Item(22, 15).Style = "Normal"
$ws.Cells.Item(22, 16).Value = ""
$ws.Cells.Item(22, 16).Style = "Normal"
$ws.Cells.Item(22, 17).Value = ""
$ws.Cells.Item(22, 17).Style = "Normal"
$ws.Cells.Item(22, 18).Value = ""
$ws.Cells.Item(22, 18).Style = "Normal"
$ws.Cells.Item(22, 19).Value = ""
$ws.Cells.Item(22, 19).Style = "Normal"
$ws.Cells.Item(22, 20).Value = ""
$ws.Cells.Item(22, 20).Style = "Normal"
$ws.Cells.Item(22, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(22, 22).Value = ""
$ws.Cells.Item(22, 22).Style = "Normal"
$ws.Cells.Item(22, 23).Value = "FLOORS-(F)"
$ws.Cells.Item(22, 24).Value = "Threshold plate-(Threshold plate)"
$ws.Cells.Item(22, 25).Value = "FLOOR BOARD NAILS FITTING 8 PES"
$ws.Cells.Item(22, 26).Value = ""
$ws.Cells.Item(22, 26).Style = "Normal"
$ws.Cells.Item(22, 27).Value = ""
$ws.Cells.Item(22, 27).Style = "Normal"
$ws.Cells.Item(22, 28).Value = ""
$ws.Cells.Item(22, 28).Style = "Normal"

# Row 26
$ws.Cells.Item(26, 2).Value = ""
$ws.Cells.Item(26, 2).Style = "Normal"
$ws.Cells.Item(26, 3).Value = ""
$ws.Cells.Item(26, 3).Style = "Normal"
$ws.Cells.Item(26, 4).Value = ""
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = ""
$ws.Cells.Item(26, 5).Style = "Normal"
$ws.Cells.Item(26, 6).Value = ""
$ws.Cells.Item(26, 6).Style = "Normal"
$ws.Cells.Item(26, 7).Value = ""
$ws.Cells.Item(26, 7).Style = "Normal"
$ws.Cells.Item(26, 8).Value = ""
$ws.Cells.Item(26, 8).Style = "Normal"
$ws.Cells.Item(26, 9).Value = ""
$ws.Cells.Item(26, 9).Style = "Normal"
$ws.Cells.Item(26, 10).Value = ""
$ws.Cells.Item(26, 10).Style = "Normal"
$ws.Cells.Item(26, 11).Value = ""
$ws.Cells.Item(26, 11).Style = "Normal"
$ws.Cells.Item(26, 12).Value = ""
$ws.Cells.Item(26, 12).Style = "Normal"
$ws.Cells.Item(26, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(26, 14).Value = ""
$ws.Cells.Item(26, 14).Style = "Normal"
$ws.Cells.Item(26, 15).Value = ""
$ws.Cells.Item(26, 15).Style = "Normal"
$ws.Cells.Item(26, 16).Value = ""
$ws.Cells.Item(26, 16).Style = "Normal"
$ws.Cells.Item(26, 17).Value = ""
$ws.Cells.Item(26, 17).Style = "Normal"
$ws.Cells.Item(26, 18).Value = ""
$ws.Cells.Item(26, 18).Style = "Normal"
$ws.Cells.Item(26, 19).Value = ""
$ws.Cells.Item(26, 19).Style = "Normal"
$ws.Cells.Item(26, 20).Value = ""
$ws.Cells.Item(26, 20).Style = "Normal"
$ws.Cells.Item(26, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(26, 22).Value = ""
$ws.Cells.Item(26, 22).Style = "Normal"
$ws.Cells.Item(26, 23).Value = "PANELS-(PANELS)"
$ws.Cells.Item(26, 24).Value = "Right side panel -(Right side panel )"
$ws.Cells.Item(26, 25).Value = "INTERIRO PANEL DIRTY"
$ws.Cells.Item(26, 26).Value = ""
$ws.Cells.Item(26, 26).Style = "Normal"
$ws.Cells.Item(26, 27).Value = ""
$ws.Cells.Item(26, 27).Style = "Normal"
$ws.Cells.Item(26, 28).Value = ""
$ws.Cells.Item(26, 28).Style = "Normal"

# Row 28
$ws.Cells.Item(28, 2).Value = ""
$ws.Cells.Item(28, 2).Style = "Normal"
$ws.Cells.Item(28, 3).Value = ""
$ws.Cells.Item(28, 3).Style = "Normal"
$ws.Cells.Item(28, 4).Value = ""
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = ""
$ws.Cells.Item(28, 5).Style = "Normal"
$ws.Cells.Item(28, 6).Value = ""
$ws.Cells.Item(28, 6).Style = "Normal"
$ws.Cells.Item(28, 7).Value = ""
$ws.Cells.Item(28, 7).Style = "Normal"
$ws.Cells.Item(28, 8).Value = ""
$ws.Cells.Item(28, 8).Style = "Normal"
$ws.Cells.Item(28, 9).Value = ""
$ws.Cells.Item(28, 9).Style = "Normal"
$ws.Cells.Item(28, 10).Value = ""
$ws.Cells.Item(28, 10).Style = "Normal"
$ws.Cells.Item(28, 11).Value = ""
$ws.Cells.Item(28, 11).Style = "Normal"
$ws.Cells.Item(28, 12).Value = ""
$ws.Cells.Item(28, 12).Style = "Normal"
$ws.Cells.Item(28, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(28, 14).Value = ""
$ws.Cells.Item(28, 14).Style = "Normal"
$ws.Cells.Item(28, 15).Value = ""
$ws.Cells.Item(28, 15).Style = "Normal"
$ws.Cells.Item(28, 16).Value = ""
$ws.Cells.Item(28, 16).Style = "Normal"
$ws.Cells.Item(28, 17).Value = ""
$ws.Cells.Item(28, 17).Style = "Normal"
$ws.Cells.Item(28, 18).Value = ""
$ws.Cells.Item(28, 18).Style = "Normal"
$ws.Cells.Item(28, 19).Value = ""
$ws.Cells.Item(28, 19).Style = "Normal"
$ws.Cells.Item(28, 20).Value = ""
$ws.Cells.Item(28, 20).Style = "Normal"
$ws.Cells.Item(28, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(28, 22).Value = ""
$ws.Cells.Item(28, 22).Style = "Normal"
$ws.Cells.Item(28, 23).Value = "PANELS-(PANELS)"
$ws.Cells.Item(28, 24).Value = "Right side panel -(Right side panel )"
$ws.Cells.Item(28, 25).Value = "INTERIOR PANEL DIRTY"
$ws.Cells.Item(28, 26).Value = ""
$ws.Cells.Item(28, 26).Style = "Normal"
$ws.Cells.Item(28, 27).Value = ""
$ws.Cells.Item(28, 27).Style = "Normal"
$ws.Cells.Item(28, 28).Value = ""
$ws.Cells.Item(28, 28).Style = "Normal"

# Row 30
$ws.Cells.Item(30, 2).Value = ""
$ws.Cells.Item(30, 2).Style = "Normal"
$ws.Cells.Item(30, 3).Value = ""
$ws.Cells.Item(30, 3).Style = "Normal"
$ws.Cells.Item(30, 4).Value = ""
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = ""
$ws.Cells.Item(30, 5).Style = "Normal"
$ws.Cells.Item(30, 6).Value = ""
$ws.Cells.Item(30, 6).Style = "Normal"
$ws.Cells.Item(30, 7).Value = ""
$ws.Cells.Item(30, 7).Style = "Normal"
$ws.Cells.Item(30, 8).Value = ""
$ws.Cells.Item(30, 8).Style = "Normal"
$ws.Cells.Item(30, 9).Value = ""
$ws.Cells.Item(30, 9).Style = "Normal"
$ws.Cells.Item(30, 10).Value = ""
$ws.Cells.Item(30, 10).Style = "Normal"
$ws.Cells.Item(30, 11).Value = ""
$ws.Cells.Item(30, 11).Style = "Normal"
$ws.Cells.Item(30, 12).Value = ""
$ws.Cells.Item(30, 12).Style = "Normal"
$ws.Cells.Item(30, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(30, 14).Value = ""
$ws.Cells.Item(30, 14).Style = "Normal"
$ws.Cells.Item(30, 15).Value = ""
$ws.Cells.Item(30, 15).Style = "Normal"
$ws.Cells.Item(30, 16).Value = ""
$ws.Cells.Item(30, 16).Style = "Normal"
$ws.Cells.Item(30, 17).Value = ""
$ws.Cells.Item(30, 17).Style = "Normal"
$ws.Cells.Item(30, 18).Value = ""
$ws.Cells.Item(30, 18).Style = "Normal"
$ws.Cells.Item(30, 19).Value = ""
$ws.Cells.Item(30, 19).Style = "Normal"
$ws.Cells.Item(30, 20).Value = ""
$ws.Cells.Item(30, 20).Style = "Normal"
$ws.Cells.Item(30, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(30, 22).Value = ""
$ws.Cells.Item(30, 22).Style = "Normal"
$ws.Cells.Item(30, 23).Value = "DOORS-(D)"
$ws.Cells.Item(30, 24).Value = "Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)"
$ws.Cells.Item(30, 25).Value = "RIGHT DOOR LOCK BAR BENT 01 PCS"
$ws.Cells.Item(30, 26).Value = ""
$ws.Cells.Item(30, 26).Style = "Normal"
$ws.Cells.Item(30, 27).Value = ""
$ws.Cells.Item(30, 27).Style = "Normal"
$ws.Cells.Item(30, 28).Value = ""
$ws.Cells.Item(30, 28).Style = "Normal"

# Row 31
$ws.Cells.Item(31, 2).Value = ""
$ws.Cells.Item(31, 2).Style = "Normal"
$ws.Cells.Item(31, 3).Value = ""
$ws.Cells.Item(31, 3).Style = "Normal"
$ws.Cells.Item(31, 4).Value = ""
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = ""
$ws.Cells.Item(31, 5).Style = "Normal"
$ws.Cells.Item(31, 6).Value = ""
$ws.Cells.Item(31, 6).Style = "Normal"
$ws.Cells.Item(31, 7).Value = ""
$ws.Cells.Item(31, 7).Style = "Normal"
$ws.Cells.Item(31, 8).Value = ""
$ws.Cells.Item(31, 8).Style = "Normal"
$ws.Cells.Item(31, 9).Value = ""
$ws.Cells.Item(31, 9).Style = "Normal"
$ws.Cells.Item(31, 10).Value = ""
$ws.Cells.Item(31, 10).Style = "Normal"
$ws.Cells.Item(31, 11).Value = ""
$ws.Cells.Item(31, 11).Style = "Normal"
$ws.Cells.Item(31, 12).Value = ""
$ws.Cells.Item(31, 12).Style = "Normal"
$ws.Cells.Item(31, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(31, 14).Value = ""
$ws.Cells.Item(31, 14).Style = "Normal"
$ws.Cells.Item(31, 15).Value = ""
$ws.Cells.Item(31, 15).Style = "Normal"
$ws.Cells.Item(31, 16).Value = ""
$ws.Cells.Item(31, 16).Style = "Normal"
$ws.Cells.Item(31, 17).Value = ""
$ws.Cells.Item(31, 17).Style = "Normal"
$ws.Cells.Item(31, 18).Value = ""
$ws.Cells.Item(31, 18).Style = "Normal"
$ws.Cells.Item(31, 19).Value = ""
$ws.Cells.Item(31, 19).Style = "Normal"
$ws.Cells.Item(31, 20).Value = ""
$ws.Cells.Item(31, 20).Style = "Normal"
$ws.Cells.Item(31, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(31, 22).Value = ""
$ws.Cells.Item(31, 22).Style = "Normal"
$ws.Cells.Item(31, 23).Value = "DOORS-(D)"
$ws.Cells.Item(31, 24).Value = "Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)"
$ws.Cells.Item(31, 25).Value = "RIGHT DOOR BOTTOM KEEPER BENT 01 PCS"
$ws.Cells.Item(31, 26).Value = ""
$ws.Cells.Item(31, 26).Style = "Normal"
$ws.Cells.Item(31, 27).Value = ""
$ws.Cells.Item(31, 27).Style = "Normal"
$ws.Cells.Item(31, 28).Value = ""
$ws.Cells.Item(31, 28).Style = "Normal"

# Row 36
$ws.Cells.Item(36, 2).Value = ""
$ws.Cells.Item(36, 2).Style = "Normal"
$ws.Cells.Item(36, 3).Value = ""
$ws.Cells.Item(36, 3).Style = "Normal"
$ws.Cells.Item(36, 4).Value = ""
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = ""
$ws.Cells.Item(36, 5).Style = "Normal"
$ws.Cells.Item(36, 6).Value = ""
$ws.Cells.Item(36, 6).Style = "Normal"
$ws.Cells.Item(36, 7).Value = ""
$ws.Cells.Item(36, 7).Style = "Normal"
$ws.Cells.Item(36, 8).Value = ""
$ws.Cells.Item(36, 8).Style = "Normal"
$ws.Cells.Item(36, 9).Value = ""
$ws.Cells.Item(36, 9).Style = "Normal"
$ws.Cells.Item(36, 10).Value = ""
$ws.Cells.Item(36, 10).Style = "Normal"
$ws.Cells.Item(36, 11).Value = ""
$ws.Cells.Item(36, 11).Style = "Normal"
$ws.Cells.Item(36, 12).Value = ""
$ws.Cells.Item(36, 12).Style = "Normal"
$ws.Cells.Item(36, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(36, 14).Value = ""
$ws.Cells.Item(36, 14).Style = "Normal"
$ws.Cells.Item(36, 15).Value = ""
$ws.Cells.Item(36, 15).Style = "Normal"
$ws.Cells.Item(36, 16).Value = ""
$ws.Cells.Item(36, 16).Style = "Normal"
$ws.Cells.Item(36, 17).Value = ""
$ws.Cells.Item(36, 17).Style = "Normal"
$ws.Cells.Item(36, 18).Value = ""
$ws.Cells.Item(36, 18).Style = "Normal"
$ws.Cells.Item(36, 19).Value = ""
$ws.Cells.Item(36, 19).Style = "Normal"
$ws.Cells.Item(36, 20).Value = ""
$ws.Cells.Item(36, 20).Style = "Normal"
$ws.Cells.Item(36, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(36, 22).Value = ""
$ws.Cells.Item(36, 22).Style = "Normal"
$ws.Cells.Item(36, 23).Value = "FLOORS-(F)"
$ws.Cells.Item(36, 24).Value = "FLOOR BOARD-(FLOOR BOARD)"
$ws.Cells.Item(36, 25).Value = "FLOOR BOARD NAILS FITTINGS AT VARIOUS PLASE"
$ws.Cells.Item(36, 26).Value = ""
$ws.Cells.Item(36, 26).Style = "Normal"
$ws.Cells.Item(36, 27).Value = ""
$ws.Cells.Item(36, 27).Style = "Normal"
$ws.Cells.Item(36, 28).Value = ""
$ws.Cells.Item(36, 28).Style = "Normal"

# Row 41
$ws.Cells.Item(41, 2).Value = ""
$ws.Cells.Item(41, 2).Style = "Normal"
$ws.Cells.Item(41, 3).Value = ""
$ws.Cells.Item(41, 3).Style = "Normal"
$ws.Cells.Item(41, 4).Value = ""
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = ""
$ws.Cells.Item(41, 5).Style = "Normal"
$ws.Cells.Item(41, 6).Value = ""
$ws.Cells.Item(41, 6).Style = "Normal"
$ws.Cells.Item(41, 7).Value = ""
$ws.Cells.Item(41, 7).Style = "Normal"
$ws.Cells.Item(41, 8).Value = ""
$ws.Cells.Item(41, 8).Style = "Normal"
$ws.Cells.Item(41, 9).Value = ""
$ws.Cells.Item(41, 9).Style = "Normal"
$ws.Cells.Item(41, 10).Value = ""
$ws.Cells.Item(41, 10).Style = "Normal"
$ws.Cells.Item(41, 11).Value = ""
$ws.Cells.Item(41, 11).Style = "Normal"
$ws.Cells.Item(41, 12).Value = ""
$ws.Cells.Item(41, 12).Style = "Normal"
$ws.Cells.Item(41, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(41, 14).Value = ""
$ws.Cells.Item(41, 14).Style = "Normal"
$ws.Cells.Item(41, 15).Value = ""
$ws.Cells.Item(41, 15).Style = "Normal"
$ws.Cells.Item(41, 16).Value = ""
$ws.Cells.Item(41, 16).Style = "Normal"
$ws.Cells.Item(41, 17).Value = ""
$ws.Cells.Item(41, 17).Style = "Normal"
$ws.Cells.Item(41, 18).Value = ""
$ws.Cells.Item(41, 18).Style = "Normal"
$ws.Cells.Item(41, 19).Value = ""
$ws.Cells.Item(41, 19).Style = "Normal"
$ws.Cells.Item(41, 20).Value = ""
$ws.Cells.Item(41, 20).Style = "Normal"
$ws.Cells.Item(41, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(41, 22).Value = ""
$ws.Cells.Item(41, 22).Style = "Normal"
$ws.Cells.Item(41, 23).Value = "DOORS-(D)"
$ws.Cells.Item(41, 24).Value = "Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)"
$ws.Cells.Item(41, 25).Value = "LEFT DOOR HINGERS PIN MISSING 01 PCS"
$ws.Cells.Item(41, 26).Value = ""
$ws.Cells.Item(41, 26).Style = "Normal"
$ws.Cells.Item(41, 27).Value = ""
$ws.Cells.Item(41, 27).Style = "Normal"
$ws.Cells.Item(41, 28).Value = ""
$ws.Cells.Item(41, 28).Style = "Normal"

# Row 45
$ws.Cells.Item(45, 2).Value = ""
$ws.Cells.Item(45, 2).Style = "Normal"
$ws.Cells.Item(45, 3).Value = ""
$ws.Cells.Item(45, 3).Style = "Normal"
$ws.Cells.Item(45, 4).Value = ""
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = ""
$ws.Cells.Item(45, 5).Style = "Normal"
$ws.Cells.Item(45, 6).Value = ""
$ws.Cells.Item(45, 6).Style = "Normal"
$ws.Cells.Item(45, 7).Value = ""
$ws.Cells.Item(45, 7).Style = "Normal"
$ws.Cells.Item(45, 8).Value = ""
$ws.Cells.Item(45, 8).Style = "Normal"
$ws.Cells.Item(45, 9).Value = ""
$ws.Cells.Item(45, 9).Style = "Normal"
$ws.Cells.Item(45, 10).Value = ""
$ws.Cells.Item(45, 10).Style = "Normal"
$ws.Cells.Item(45, 11).Value = ""
$ws.Cells.Item(45, 11).Style = "Normal"
$ws.Cells.Item(45, 12).Value = ""
$ws.Cells.Item(45, 12).Style = "Normal"
$ws.Cells.Item(45, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(45, 14).Value = ""
$ws.Cells.Item(45, 14).Style = "Normal"
$ws.Cells.Item(45, 15).Value = ""
$ws.Cells.Item(45, 15).Style = "Normal"
$ws.Cells.Item(45, 16).Value = ""
$ws.Cells.Item(45, 16).Style = "Normal"
$ws.Cells.Item(45, 17).Value = ""
$ws.Cells.Item(45, 17).Style = "Normal"
$ws.Cells.Item(45, 18).Value = ""
$ws.Cells.Item(45, 18).Style = "Normal"
$ws.Cells.Item(45, 19).Value = ""
$ws.Cells.Item(45, 19).Style = "Normal"
$ws.Cells.Item(45, 20).Value = ""
$ws.Cells.Item(45, 20).Style = "Normal"
$ws.Cells.Item(45, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(45, 22).Value = ""
$ws.Cells.Item(45, 22).Style = "Normal"
$ws.Cells.Item(45, 23).Value = "FLOORS-(F)"
$ws.Cells.Item(45, 24).Value = "Threshold plate-(Threshold plate)"
$ws.Cells.Item(45, 25).Value = "FLOOR BOARD NAILS FITTING 03 PLS"
$ws.Cells.Item(45, 26).Value = ""
$ws.Cells.Item(45, 26).Style = "Normal"
$ws.Cells.Item(45, 27).Value = ""
$ws.Cells.Item(45, 27).Style = "Normal"
$ws.Cells.Item(45, 28).Value = ""
$ws.Cells.Item(45, 28).Style = "Normal"

# Row 47
$ws.Cells.Item(47, 2).Value = ""
$ws.Cells.Item(47, 2).Style = "Normal"
$ws.Cells.Item(47, 3).Value = ""
$ws.Cells.Item(47, 3).Style = "Normal"
$ws.Cells.Item(47, 4).Value = ""
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = ""
$ws.Cells.Item(47, 5).Style = "Normal"
$ws.Cells.Item(47, 6).Value = ""
$ws.Cells.Item(47, 6).Style = "Normal"
$ws.Cells.Item(47, 7).Value = ""
$ws.Cells.Item(47, 7).Style = "Normal"
$ws.Cells.Item(47, 8).Value = ""
$ws.Cells.Item(47, 8).Style = "Normal"
$ws.Cells.Item(47, 9).Value = ""
$ws.Cells.Item(47, 9).Style = "Normal"
$ws.Cells.Item(47, 10).Value = ""
$ws.Cells.Item(47, 10).Style = "Normal"
$ws.Cells.Item(47, 11).Value = ""
$ws.Cells.Item(47, 11).Style = "Normal"
$ws.Cells.Item(47, 12).Value = ""
$ws.Cells.Item(47, 12).Style = "Normal"
$ws.Cells.Item(47, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(47, 14).Value = ""
$ws.Cells.Item(47, 14).Style = "Normal"
$ws.Cells.Item(47, 15).Value = ""
$ws.Cells.Item(47, 15).Style = "Normal"
$ws.Cells.Item(47, 16).Value = ""
$ws.Cells.Item(47, 16).Style = "Normal"
$ws.Cells.Item(47, 17).Value = ""
$ws.Cells.Item(47, 17).Style = "Normal"
$ws.Cells.Item(47, 18).Value = ""
$ws.Cells.Item(47, 18).Style = "Normal"
$ws.Cells.Item(47, 19).Value = ""
$ws.Cells.Item(47, 19).Style = "Normal"
$ws.Cells.Item(47, 20).Value = ""
$ws.Cells.Item(47, 20).Style = "Normal"
$ws.Cells.Item(47, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(47, 22).Value = ""
$ws.Cells.Item(47, 22).Style = "Normal"
$ws.Cells.Item(47, 23).Value = "PANELS-(PANELS)"
$ws.Cells.Item(47, 24).Value = "Right side panel -(Right side panel )"
$ws.Cells.Item(47, 25).Value = "INTERIOR PANEL DIRTY"
$ws.Cells.Item(47, 26).Value = ""
$ws.Cells.Item(47, 26).Style = "Normal"
$ws.Cells.Item(47, 27).Value = ""
$ws.Cells.Item(47, 27).Style = "Normal"
$ws.Cells.Item(47, 28).Value = ""
$ws.Cells.Item(47, 28).Style = "Normal"

# Row 49
$ws.Cells.Item(49, 2).Value = ""
$ws.Cells.Item(49, 2).Style = "Normal"
$ws.Cells.Item(49, 3).Value = ""
$ws.Cells.Item(49, 3).Style = "Normal"
$ws.Cells.Item(49, 4).Value = ""
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = ""
$ws.Cells.Item(49, 5).Style = "Normal"
$ws.Cells.Item(49, 6).Value = ""
$ws.Cells.Item(49, 6).Style = "Normal"
$ws.Cells.Item(49, 7).Value = ""
$ws.Cells.Item(49, 7).Style = "Normal"
$ws.Cells.Item(49, 8).Value = ""
$ws.Cells.Item(49, 8).Style = "Normal"
$ws.Cells.Item(49, 9).Value = ""
$ws.Cells.Item(49, 9).Style = "Normal"
$ws.Cells.Item(49, 10).Value = ""
$ws.Cells.Item(49, 10).Style = "Normal"
$ws.Cells.Item(49, 11).Value = ""
$ws.Cells.Item(49, 11).Style = "Normal"
$ws.Cells.Item(49, 12).Value = ""
$ws.Cells.Item(49, 12).Style = "Normal"
$ws.Cells.Item(49, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(49, 14).Value = ""
$ws.Cells.Item(49, 14).Style = "Normal"
$ws.Cells.Item(49, 15).Value = ""
$ws.Cells.Item(49, 15).Style = "Normal"
$ws.Cells.Item(49, 16).Value = ""
$ws.Cells.Item(49, 16).Style = "Normal"
$ws.Cells.Item(49, 17).Value = ""
$ws.Cells.Item(49, 17).Style = "Normal"
$ws.Cells.Item(49, 18).Value = ""
$ws.Cells.Item(49, 18).Style = "Normal"
$ws.Cells.Item(49, 19).Value = ""
$ws.Cells.Item(49, 19).Style = "Normal"
$ws.Cells.Item(49, 20).Value = ""
$ws.Cells.Item(49, 20).Style = "Normal"
$ws.Cells.Item(49, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(49, 22).Value = ""
$ws.Cells.Item(49, 22).Style = "Normal"
$ws.Cells.Item(49, 23).Value = "PANELS-(PANELS)"
$ws.Cells.Item(49, 24).Value = "Right side panel -(Right side panel )"
$ws.Cells.Item(49, 25).Value = "INTERIOR PANEL INK"
$ws.Cells.Item(49, 26).Value = ""
$ws.Cells.Item(49, 26).Style = "Normal"
$ws.Cells.Item(49, 27).Value = ""
$ws.Cells.Item(49, 27).Style = "Normal"
$ws.Cells.Item(49, 28).Value = ""
$ws.Cells.Item(49, 28).Style = "Normal"

# Row 51
$ws.Cells.Item(51, 2).Value = ""
$ws.Cells.Item(51, 2).Style = "Normal"
$ws.Cells.Item(51, 3).Value = ""
$ws.Cells.Item(51, 3).Style = "Normal"
$ws.Cells.Item(51, 4).Value = ""
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = ""
$ws.Cells.Item(51, 5).Style = "Normal"
$ws.Cells.Item(51, 6).Value = ""
$ws.Cells.Item(51, 6).Style = "Normal"
$ws.Cells.Item(51, 7).Value = ""
$ws.Cells.Item(51, 7).Style = "Normal"
$ws.Cells.Item(51, 8).Value = ""
$ws.Cells.Item(51, 8).Style = "Normal"
$ws.Cells.Item(51, 9).Value = ""
$ws.Cells.Item(51, 9).Style = "Normal"
$ws.Cells.Item(51, 10).Value = ""
$ws.Cells.Item(51, 10).Style = "Normal"
$ws.Cells.Item(51, 11).Value = ""
$ws.Cells.Item(51, 11).Style = "Normal"
$ws.Cells.Item(51, 12).Value = ""
$ws.Cells.Item(51, 12).Style = "Normal"
$ws.Cells.Item(51, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(51, 14).Value = ""
$ws.Cells.Item(51, 14).Style = "Normal"
$ws.Cells.Item(51, 15).Value = ""
$ws.Cells.Item(51, 15).Style = "Normal"
$ws.Cells.Item(51, 16).Value = ""
$ws.Cells.Item(51, 16).Style = "Normal"
$ws.Cells.Item(51, 17).Value = ""
$ws.Cells.Item(51, 17).Style = "Normal"
$ws.Cells.Item(51, 18).Value = ""
$ws.Cells.Item(51, 18).Style = "Normal"
$ws.Cells.Item(51, 19).Value = ""
$ws.Cells.Item(51, 19).Style = "Normal"
$ws.Cells.Item(51, 20).Value = ""
$ws.Cells.Item(51, 20).Style = "Normal"
$ws.Cells.Item(51, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(51, 22).Value = ""
$ws.Cells.Item(51, 22).Style = "Normal"
$ws.Cells.Item(51, 23).Value = "DOORS-(D)"
$ws.Cells.Item(51, 24).Value = "Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)"
$ws.Cells.Item(51, 25).Value = "LEFT DOOR HINGS BENT 1 PCS (DOOR UNLOCK)"
$ws.Cells.Item(51, 26).Value = ""
$ws.Cells.Item(51, 26).Style = "Normal"
$ws.Cells.Item(51, 27).Value = ""
$ws.Cells.Item(51, 27).Style = "Normal"
$ws.Cells.Item(51, 28).Value = ""
$ws.Cells.Item(51, 28).Style = "Normal"

# Row 53
$ws.Cells.Item(53, 2).Value = ""
$ws.Cells.Item(53, 2).Style = "Normal"
$ws.Cells.Item(53, 3).Value = ""
$ws.Cells.Item(53, 3).Style = "Normal"
$ws.Cells.Item(53, 4).Value = ""
$ws.Cells.Item(53, 4).Style = "Normal"
$ws.Cells.Item(53, 5).Value = ""
$ws.Cells.Item(53, 5).Style = "Normal"
$ws.Cells.Item(53, 6).Value = ""
$ws.Cells.Item(53, 6).Style = "Normal"
$ws.Cells.Item(53, 7).Value = ""
$ws.Cells.Item(53, 7).Style = "Normal"
$ws.Cells.Item(53, 8).Value = ""
$ws.Cells.Item(53, 8).Style = "Normal"
$ws.Cells.Item(53, 9).Value = ""
$ws.Cells.Item(53, 9).Style = "Normal"
$ws.Cells.Item(53, 10).Value = ""
$ws.Cells.Item(53, 10).Style = "Normal"
$ws.Cells.Item(53, 11).Value = ""
$ws.Cells.Item(53, 11).Style = "Normal"
$ws.Cells.Item(53, 12).Value = ""
$ws.Cells.Item(53, 12).Style = "Normal"
$ws.Cells.Item(53, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(53, 14).Value = ""
$ws.Cells.Item(53, 14).Style = "Normal"
$ws.Cells.Item(53, 15).Value = ""
$ws.Cells.Item(53, 15).Style = "Normal"
$ws.Cells.Item(53, 16).Value = ""
$ws.Cells.Item(53, 16).Style = "Normal"
$ws.Cells.Item(53, 17).Value = ""
$ws.Cells.Item(53, 17).Style = "Normal"
$ws.Cells.Item(53, 18).Value = ""
$ws.Cells.Item(53, 18).Style = "Normal"
$ws.Cells.Item(53, 19).Value = ""
$ws.Cells.Item(53, 19).Style = "Normal"
$ws.Cells.Item(53, 20).Value = ""
$ws.Cells.Item(53, 20).Style = "Normal"
$ws.Cells.Item(53, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(53, 22).Value = ""
$ws.Cells.Item(53, 22).Style = "Normal"
$ws.Cells.Item(53, 23).Value = "PANELS-(PANELS)"
$ws.Cells.Item(53, 24).Value = "Right side panel -(Right side panel )"
$ws.Cells.Item(53, 25).Value = "INTERIOR PANEL DIRTY"
$ws.Cells.Item(53, 26).Value = ""
$ws.Cells.Item(53, 26).Style = "Normal"
$ws.Cells.Item(53, 27).Value = ""
$ws.Cells.Item(53, 27).Style = "Normal"
$ws.Cells.Item(53, 28).Value = ""
$ws.Cells.Item(53, 28).Style = "Normal"

# Row 57
$ws.Cells.Item(57, 2).Value = ""
$ws.Cells.Item(57, 2).Style = "Normal"
$ws.Cells.Item(57, 3).Value = ""
$ws.Cells.Item(57, 3).Style = "Normal"
$ws.Cells.Item(57, 4).Value = ""
$ws.Cells.Item(57, 4).Style = "Normal"
$ws.Cells.Item(57, 5).Value = ""
$ws.Cells.Item(57, 5).Style = "Normal"
$ws.Cells.Item(57, 6).Value = ""
$ws.Cells.Item(57, 6).Style = "Normal"
$ws.Cells.Item(57, 7).Value = ""
$ws.Cells.Item(57, 7).Style = "Normal"
$ws.Cells.Item(57, 8).Value = ""
$ws.Cells.Item(57, 8).Style = "Normal"
$ws.Cells.Item(57, 9).Value = ""
$ws.Cells.Item(57, 9).Style = "Normal"
$ws.Cells.Item(57, 10).Value = ""
$ws.Cells.Item(57, 10).Style = "Normal"
$ws.Cells.Item(57, 11).Value = ""
$ws.Cells.Item(57, 11).Style = "Normal"
$ws.Cells.Item(57, 12).Value = ""
$ws.Cells.Item(57, 12).Style = "Normal"
$ws.Cells.Item(57, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(57, 14).Value = ""
$ws.Cells.Item(57, 14).Style = "Normal"
$ws.Cells.Item(57, 15).Value = ""
$ws.Cells.Item(57, 15).Style = "Normal"
$ws.Cells.Item(57, 16).Value = ""
$ws.Cells.Item(57, 16).Style = "Normal"
$ws.Cells.Item(57, 17).Value = ""
$ws.Cells.Item(57, 17).Style = "Normal"
$ws.Cells.Item(57, 18).Value = ""
$ws.Cells.Item(57, 18).Style = "Normal"
$ws.Cells.Item(57, 19).Value = ""
$ws.Cells.Item(57, 19).Style = "Normal"
$ws.Cells.Item(57, 20).Value = ""
$ws.Cells.Item(57, 20).Style = "Normal"
$ws.Cells.Item(57, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(57, 22).Value = ""
$ws.Cells.Item(57, 22).Style = "Normal"
$ws.Cells.Item(57, 23).Value = "PANELS-(PANELS)"
$ws.Cells.Item(57, 24).Value = "Right side panel -(Right side panel )"
$ws.Cells.Item(57, 25).Value = "INTERIOR PANEL DIRTY"
$ws.Cells.Item(57, 26).Value = ""
$ws.Cells.Item(57, 26).Style = "Normal"
$ws.Cells.Item(57, 27).Value = ""
$ws.Cells.Item(57, 27).Style = "Normal"
$ws.Cells.Item(57, 28).Value = ""
$ws.Cells.Item(57, 28).Style = "Normal"

# Row 58
$ws.Cells.Item(58, 2).Value = ""
$ws.Cells.Item(58, 2).Style = "Normal"
$ws.Cells.Item(58, 3).Value = ""
$ws.Cells.Item(58, 3).Style = "Normal"
$ws.Cells.Item(58, 4).Value = ""
$ws.Cells.Item(58, 4).Style = "Normal"
$ws.Cells.Item(58, 5).Value = ""
$ws.Cells.Item(58, 5).Style = "Normal"
$ws.Cells.Item(58, 6).Value = ""
$ws.Cells.Item(58, 6).Style = "Normal"
$ws.Cells.Item(58, 7).Value = ""
$ws.Cells.Item(58, 7).Style = "Normal"
$ws.Cells.Item(58, 8).Value = ""
$ws.Cells.Item(58, 8).Style = "Normal"
$ws.Cells.Item(58, 9).Value = ""
$ws.Cells.Item(58, 9).Style = "Normal"
$ws.Cells.Item(58, 10).Value = ""
$ws.Cells.Item(58, 10).Style = "Normal"
$ws.Cells.Item(58, 11).Value = ""
$ws.Cells.Item(58, 11).Style = "Normal"
$ws.Cells.Item(58, 12).Value = ""
$ws.Cells.Item(58, 12).Style = "Normal"
$ws.Cells.Item(58, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(58, 14).Value = ""
$ws.Cells.Item(58, 14).Style = "Normal"
$ws.Cells.Item(58, 15).Value = ""
$ws.Cells.Item(58, 15).Style = "Normal"
$ws.Cells.Item(58, 16).Value = ""
$ws.Cells.Item(58, 16).Style = "Normal"
$ws.Cells.Item(58, 17).Value = ""
$ws.Cells.Item(58, 17).Style = "Normal"
$ws.Cells.Item(58, 18).Value = ""
$ws.Cells.Item(58, 18).Style = "Normal"
$ws.Cells.Item(58, 19).Value = ""
$ws.Cells.Item(58, 19).Style = "Normal"
$ws.Cells.Item(58, 20).Value = ""
$ws.Cells.Item(58, 20).Style = "Normal"
$ws.Cells.Item(58, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(58, 22).Value = ""
$ws.Cells.Item(58, 22).Style = "Normal"
$ws.Cells.Item(58, 23).Value = "DOORS-(D)"
$ws.Cells.Item(58, 24).Value = "Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)"
$ws.Cells.Item(58, 25).Value = "RIGHT DOOR BOTTOM GASKET CUT 3"""
$ws.Cells.Item(58, 26).Value = ""
$ws.Cells.Item(58, 26).Style = "Normal"
$ws.Cells.Item(58, 27).Value = ""
$ws.Cells.Item(58, 27).Style = "Normal"
$ws.Cells.Item(58, 28).Value = ""
$ws.Cells.Item(58, 28).Style = "Normal"

# Row 60
$ws.Cells.Item(60, 2).Value = ""
$ws.Cells.Item(60, 2).Style = "Normal"
$ws.Cells.Item(60, 3).Value = ""
$ws.Cells.Item(60, 3).Style = "Normal"
$ws.Cells.Item(60, 4).Value = ""
$ws.Cells.Item(60, 4).Style = "Normal"
$ws.Cells.Item(60, 5).Value = ""
$ws.Cells.Item(60, 5).Style = "Normal"
$ws.Cells.Item(60, 6).Value = ""
$ws.Cells.Item(60, 6).Style = "Normal"
$ws.Cells.Item(60, 7).Value = ""
$ws.Cells.Item(60, 7).Style = "Normal"
$ws.Cells.Item(60, 8).Value = ""
$ws.Cells.Item(60, 8).Style = "Normal"
$ws.Cells.Item(60, 9).Value = ""
$ws.Cells.Item(60, 9).Style = "Normal"
$ws.Cells.Item(60, 10).Value = ""
$ws.Cells.Item(60, 10).Style = "Normal"
$ws.Cells.Item(60, 11).Value = ""
$ws.Cells.Item(60, 11).Style = "Normal"
$ws.Cells.Item(60, 12).Value = ""
$ws.Cells.Item(60, 12).Style = "Normal"
$ws.Cells.Item(60, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(60, 14).Value = ""
$ws.Cells.Item(60, 14).Style = "Normal"
$ws.Cells.Item(60, 15).Value = ""
$ws.Cells.Item(60, 15).Style = "Normal"
$ws.Cells.Item(60, 16).Value = ""
$ws.Cells.Item(60, 16).Style = "Normal"
$ws.Cells.Item(60, 17).Value = ""
$ws.Cells.Item(60, 17).Style = "Normal"
$ws.Cells.Item(60, 18).Value = ""
$ws.Cells.Item(60, 18).Style = "Normal"
$ws.Cells.Item(60, 19).Value = ""
$ws.Cells.Item(60, 19).Style = "Normal"
$ws.Cells.Item(60, 20).Value = ""
$ws.Cells.Item(60, 20).Style = "Normal"
$ws.Cells.Item(60, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(60, 22).Value = ""
$ws.Cells.Item(60, 22).Style = "Normal"
$ws.Cells.Item(60, 23).Value = "DOORS-(D)"
$ws.Cells.Item(60, 24).Value = "Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)"
$ws.Cells.Item(60, 25).Value = "LEFT DOOR TOP GASKET CUT 2"""
$ws.Cells.Item(60, 26).Value = ""
$ws.Cells.Item(60, 26).Style = "Normal"
$ws.Cells.Item(60, 27).Value = ""
$ws.Cells.Item(60, 27).Style = "Normal"
$ws.Cells.Item(60, 28).Value = ""
$ws.Cells.Item(60, 28).Style = "Normal"

# Row 61
$ws.Cells.Item(61, 2).Value = ""
$ws.Cells.Item(61, 2).Style = "Normal"
$ws.Cells.Item(61, 3).Value = ""
$ws.Cells.Item(61, 3).Style = "Normal"
$ws.Cells.Item(61, 4).Value = ""
$ws.Cells.Item(61, 4).Style = "Normal"
$ws.Cells.Item(61, 5).Value = ""
$ws.Cells.Item(61, 5).Style = "Normal"
$ws.Cells.Item(61, 6).Value = ""
$ws.Cells.Item(61, 6).Style = "Normal"
$ws.Cells.Item(61, 7).Value = ""
$ws.Cells.Item(61, 7).Style = "Normal"
$ws.Cells.Item(61, 8).Value = ""
$ws.Cells.Item(61, 8).Style = "Normal"
$ws.Cells.Item(61, 9).Value = ""
$ws.Cells.Item(61, 9).Style = "Normal"
$ws.Cells.Item(61, 10).Value = ""
$ws.Cells.Item(61, 10).Style = "Normal"
$ws.Cells.Item(61, 11).Value = ""
$ws.Cells.Item(61, 11).Style = "Normal"
$ws.Cells.Item(61, 12).Value = ""
$ws.Cells.Item(61, 12).Style = "Normal"
$ws.Cells.Item(61, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(61, 14).Value = ""
$ws.Cells.Item(61, 14).Style = "Normal"
$ws.Cells.Item(61, 15).Value = ""
$ws.Cells.Item(61, 15).Style = "Normal"
$ws.Cells.Item(61, 16).Value = ""
$ws.Cells.Item(61, 16).Style = "Normal"
$ws.Cells.Item(61, 17).Value = ""
$ws.Cells.Item(61, 17).Style = "Normal"
$ws.Cells.Item(61, 18).Value = ""
$ws.Cells.Item(61, 18).Style = "Normal"
$ws.Cells.Item(61, 19).Value = ""
$ws.Cells.Item(61, 19).Style = "Normal"
$ws.Cells.Item(61, 20).Value = ""
$ws.Cells.Item(61, 20).Style = "Normal"
$ws.Cells.Item(61, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(61, 22).Value = ""
$ws.Cells.Item(61, 22).Style = "Normal"
$ws.Cells.Item(61, 23).Value = "DOORS-(D)"
$ws.Cells.Item(61, 24).Value = "Door stiffeners hinges side edge.-(Door stiffeners hinges side edge.)"
$ws.Cells.Item(61, 25).Value = "LEFT DOOR LOCKBAR BOTTOM KEEPER BENT 1 PCS WITH LOCKBAR"
$ws.Cells.Item(61, 26).Value = ""
$ws.Cells.Item(61, 26).Style = "Normal"
$ws.Cells.Item(61, 27).Value = ""
$ws.Cells.Item(61, 27).Style = "Normal"
$ws.Cells.Item(61, 28).Value = ""
$ws.Cells.Item(61, 28).Style = "Normal"

# Row 63
$ws.Cells.Item(63, 2).Value = ""
$ws.Cells.Item(63, 2).Style = "Normal"
$ws.Cells.Item(63, 3).Value = ""
$ws.Cells.Item(63, 3).Style = "Normal"
$ws.Cells.Item(63, 4).Value = ""
$ws.Cells.Item(63, 4).Style = "Normal"
$ws.Cells.Item(63, 5).Value = ""
$ws.Cells.Item(63, 5).Style = "Normal"
$ws.Cells.Item(63, 6).Value = ""
$ws.Cells.Item(63, 6).Style = "Normal"
$ws.Cells.Item(63, 7).Value = ""
$ws.Cells.Item(63, 7).Style = "Normal"
$ws.Cells.Item(63, 8).Value = ""
$ws.Cells.Item(63, 8).Style = "Normal"
$ws.Cells.Item(63, 9).Value = ""
$ws.Cells.Item(63, 9).Style = "Normal"
$ws.Cells.Item(63, 10).Value = ""
$ws.Cells.Item(63, 10).Style = "Normal"
$ws.Cells.Item(63, 11).Value = ""
$ws.Cells.Item(63, 11).Style = "Normal"
$ws.Cells.Item(63, 12).Value = ""
$ws.Cells.Item(63, 12).Style = "Normal"
$ws.Cells.Item(63, 13).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(63, 14).Value = ""
$ws.Cells.Item(63, 14).Style = "Normal"
$ws.Cells.Item(63, 15).Value = ""
$ws.Cells.Item(63, 15).Style = "Normal"
$ws.Cells.Item(63, 16).Value = ""
$ws.Cells.Item(63, 16).Style = "Normal"
$ws.Cells.Item(63, 17).Value = ""
$ws.Cells.Item(63, 17).Style = "Normal"
$ws.Cells.Item(63, 18).Value = ""
$ws.Cells.Item(63, 18).Style = "Normal"
$ws.Cells.Item(63, 19).Value = ""
$ws.Cells.Item(63, 19).Style = "Normal"
$ws.Cells.Item(63, 20).Value = ""
$ws.Cells.Item(63, 20).Style = "Normal"
$ws.Cells.Item(63, 21).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(63, 22).Value = ""
$ws.Cells.Item(63, 22).Style = "Normal"
$ws.Cells.Item(63, 23).Value = "PANELS-(PANELS)"
$ws.Cells.Item(63, 24).Value = "Right side panel -(Right side panel )"
$ws.Cells.Item(63, 25).Value = "INTERIOR PANEL DIRTY"
$ws.Cells.Item(63, 26).Value = ""
$ws.Cells.Item(63, 26).Style = "Normal"
$ws.Cells.Item(63, 27).Value = ""
$ws.Cells.Item(63, 27).Style = "Normal"
$ws.Cells.Item(63, 28).Value = ""
$ws.Cells.Item(63, 28).Style = "Normal"

# Column-width follow-up: column B no longer holds the giant combined-array
# text (so it can shrink), while columns X and Y now hold the longer
# "Damage Component" / "Remark" text that used to live inside column B,
# so they need to grow to fit it.
$ws.Columns.Item(2).ColumnWidth = 7.88
$ws.Columns.Item(24).ColumnWidth = 46.38
$ws.Columns.Item(25).ColumnWidth = 57.38
